# "letter A for Green"
# Adds a new "Sheet3" after the existing "Sheet2", makes it the active/selected
# sheet, and fills B2:P16 with a 0/1 grid that draws the letter "A"
# (mirrors the style already used on the other sheets of this workbook).

$wb = $excel.ActiveWorkbook

# Insert the new worksheet right after Sheet2 so it becomes Sheet3 and the
# last (active) tab, exactly like the diff shows.
$sheet2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet2)

# 0/1 bitmap for the letter "A", rows 2-16, columns B(2)-P(16).
$grid = @(
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,1,1,0,0,0,0,0,0),
    @(0,0,0,0,0,0,1,1,1,1,0,0,0,0,0),
    @(0,0,0,0,0,1,1,1,1,1,1,0,0,0,0),
    @(0,0,0,0,1,1,1,0,0,1,1,1,0,0,0),
    @(0,0,0,0,1,1,0,0,0,0,1,1,0,0,0),
    @(0,0,0,0,1,1,0,0,0,0,1,1,0,0,0),
    @(0,0,0,0,1,1,1,1,1,1,1,1,0,0,0),
    @(0,0,0,0,1,1,1,1,1,1,1,1,0,0,0),
    @(0,0,0,0,1,1,0,0,0,0,1,1,0,0,0),
    @(0,0,0,0,1,1,0,0,0,0,1,1,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
)

for ($i = 0; $i -lt $grid.Count; $i++) {
    $row = $grid[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $c = $j + 2
        $ws3.Cells.Item($r, $c).Value = $row[$j]
    }
}

# Select B2:P16 on the new sheet so it is the recorded selection/tabSelected view.
$ws3.Range("B2:P16").Select()
